$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------------
# The sheet used to show a 3-year comparison table (1989 / 2002 / 2014)
# with an extra subtitle line ("(according to the population census
# data)"). The export was fixed/regenerated so the sheet now only shows
# the single current value (2014) and the subtitle row was dropped.
# -----------------------------------------------------------------------

# Drop the subtitle text in A2 ("(according to the population census
# data)") entirely -- row 2 stays as a blank spacer row.
$ws.Range("A2").Clear()

# Remove the blank spacer row that used to sit between the title block
# and the "(sq. km)" label (old row 3) -- this shifts the "(sq. km)"
# label, the year headers and the data row up by one.
$ws.Rows(3).Delete()

# Only the most recent year (2014) is kept; drop the 1989 and 2002
# columns so the 2014 column becomes column B.
$ws.Range("B:C").EntireColumn.Delete()

# Match the taller row height used by the regenerated export.
$ws.Range("A1:B8").RowHeight = 20.1

Write-Host "done"
